$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.780.09"
$ws.Range("E2").Value = "  -3.36%  "
$ws.Range("D3").Value = "2.900.73"
$ws.Range("E3").Value = "  -4.47%  "
$ws.Range("E4").Value = "  +0.12%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "586.24"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.88%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "145.27"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -4.51%  "
$ws.Range("E7").Value = "  +0.12%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.500"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -3.76%  "
$ws.Range("D9").Value = "2.899.32"
$ws.Range("E9").Value = "  -4.43%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "6.66"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +3.86%  "
$ws.Range("E11").Value = "  -5.86%  "
$ws.Range("E12").Value = "  -3.56%  "
$ws.Range("E13").Value = "  -5.08%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "33.38"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -4.45%  "
$ws.Range("E15").Value = "  +0.49%  "
$ws.Range("D16").Value = "3.382.87"
$ws.Range("E16").Value = "  -4.15%  "
$ws.Range("D17").Value = "60.675.25"
$ws.Range("E17").Value = "  -3.25%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "6.71"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -4.65%  "
$ws.Range("D19").Value = "2.900.48"
$ws.Range("E19").Value = "  -4.43%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "423.79"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -6.36%  "
$ws.Range("E21").Value = "  -5.16%  "
$ws.Range("E22").Value = "  -3.68%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "7.05"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -5.98%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "79.82"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -3.24%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "10.97"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -3.29%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "11.79"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -3.71%  "
$ws.Range("E28").Value = "  +0.02%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("E30").Value = "  -3.59%  "
$ws.Range("E31").Value = "  -4.13%  "
$ws.Range("E32").Value = "  -0.34%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "26.21"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -5.29%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.105"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -4.64%  "
$ws.Range("D35").Value = "0.0₃0831"
$ws.Range("E35").Value = "  -4.31%  "
$ws.Range("E36").Value = "  -3.08%  "
$ws.Range("E37").Value = "  -5.22%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "49.33"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -2.10%  "
$ws.Range("E39").Value = "  -5.74%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "2.01"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -4.93%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.124"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.58%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "8.63"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -4.94%  "
$ws.Range("E43").Value = "  -2.31%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "41.06"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.41%  "
$ws.Range("E45").Value = "  -3.01%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "371.89"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -5.50%  "
$ws.Range("D47").Value = "2.660.04"
$ws.Range("E47").Value = "  -3.22%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "132.89"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +0.59%  "
$ws.Range("E49").Value = "  +0.00%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "25.07"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +4.22%  "
$ws.Range("E51").Value = "  -2.32%  "
